# Auto-generated edit script: applies the Excalibur_Profits.xlsx value
# updates (scheduled-runner price refresh) described by the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 542.8570999999999
$ws.Range("J12").Value = 733.3333
$ws.Range("L12").Value = 733.3333
$ws.Range("N12").Value = -1073.3333
$ws.Range("H18").Value = 1299.625
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null
$ws.Range("H28").Value = 145835.72
$ws.Range("I28").Value = 202380
$ws.Range("J28").Value = 4475
$ws.Range("K28").Value = 202380
$ws.Range("L28").Value = 4475
$ws.Range("M28").Value = -201895
$ws.Range("N28").Value = -5445
$ws.Range("H32").Value = 5049.375
$ws.Range("J32").Value = 5039.4
$ws.Range("L32").Value = 5039.4
$ws.Range("N32").Value = -5691.4
$ws.Range("H40").Value = 6342.091
$ws.Range("I40").Value = 7681.2856
$ws.Range("J40").Value = 3998.5
$ws.Range("K40").Value = 7681.2856
$ws.Range("L40").Value = 3998.5
$ws.Range("M40").Value = -7506.2856
$ws.Range("N40").Value = -4348.5
$ws.Range("H51").Value = 12733.333
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H69").Value = 12645.044
$ws.Range("I69").Value = 5599.5
$ws.Range("J69").Value = 14128.315
$ws.Range("K69").Value = 16798.5
$ws.Range("L69").Value = 42384.945
$ws.Range("M69").Value = -15924.5
$ws.Range("N69").Value = -44132.945
$ws.Range("H70").Value = 1265.2667
$ws.Range("I70").Value = 1152.6364
$ws.Range("J70").Value = 1575
$ws.Range("K70").Value = 3457.9092
$ws.Range("L70").Value = 4725
$ws.Range("M70").Value = -3187.9092
$ws.Range("N70").Value = -5265
$ws.Range("H72").Value = 12645.044
$ws.Range("I72").Value = 5599.5
$ws.Range("J72").Value = 14128.315
$ws.Range("K72").Value = 50395.5
$ws.Range("L72").Value = 127154.835
$ws.Range("M72").Value = -46027.5
$ws.Range("N72").Value = -135890.835
$ws.Range("H73").Value = 1265.2667
$ws.Range("I73").Value = 1152.6364
$ws.Range("J73").Value = 1575
$ws.Range("K73").Value = 3457.9092
$ws.Range("L73").Value = 4725
$ws.Range("M73").Value = -2521.9092
$ws.Range("N73").Value = -6597
$ws.Range("H92").Value = 674.2692
$ws.Range("I92").Value = 551.4545000000001
$ws.Range("J92").Value = 1349.75
$ws.Range("K92").Value = 551.4545000000001
$ws.Range("L92").Value = 1349.75
$ws.Range("M92").Value = 696.5454999999999
$ws.Range("N92").Value = -3845.75
$ws.Range("H132").Value = 7779.1313
$ws.Range("I132").Value = 7854.243
$ws.Range("K132").Value = 23562.729
$ws.Range("M132").Value = -21032.729
$ws.Range("H138").Value = 3456.3635
$ws.Range("I138").Value = 1485.9
$ws.Range("J138").Value = 4035.9119
$ws.Range("K138").Value = 4457.700000000001
$ws.Range("L138").Value = 12107.7357
$ws.Range("M138").Value = 682.2999999999993
$ws.Range("N138").Value = -22387.7357
$ws.Range("H141").Value = 1114.8823
$ws.Range("I141").Value = 1146.75
$ws.Range("K141").Value = 3440.25
$ws.Range("M141").Value = 1739.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 828.8936
$ws.Range("I32").Value = 576.25
$ws.Range("K32").Value = 576.25
$ws.Range("M32").Value = -289.25
$ws.Range("H45").Value = 3637.3333
$ws.Range("I45").Value = 3956.3333
$ws.Range("K45").Value = 3956.3333
$ws.Range("M45").Value = -3579.3333
$ws.Range("H61").Value = 4011.4
$ws.Range("I61").Value = 4011.4
$ws.Range("K61").Value = 4011.4
$ws.Range("M61").Value = -3799.4
$ws.Range("H122").Value = 3296.0344
$ws.Range("I122").Value = 2689.8
$ws.Range("J122").Value = 4643.222
$ws.Range("K122").Value = 8069.400000000001
$ws.Range("L122").Value = 13929.666
$ws.Range("M122").Value = -5619.400000000001
$ws.Range("N122").Value = -18829.666
$ws.Range("H132").Value = 3984.5
$ws.Range("I132").Value = 3632.95
$ws.Range("K132").Value = 10898.85
$ws.Range("M132").Value = -8368.849999999999
$ws.Range("H136").Value = 4011.4
$ws.Range("I136").Value = 4011.4
$ws.Range("K136").Value = 12034.2
$ws.Range("M136").Value = -9484.200000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 692.75
$ws.Range("I22").Value = 824.1667
$ws.Range("K22").Value = 824.1667
$ws.Range("M22").Value = -651.1667
$ws.Range("H86").Value = 3227.0588
$ws.Range("I86").Value = 3558.6667
$ws.Range("J86").Value = 2431.2
$ws.Range("K86").Value = 3558.6667
$ws.Range("L86").Value = 2431.2
$ws.Range("M86").Value = -2435.6667
$ws.Range("N86").Value = -4677.2
$ws.Range("H89").Value = 3227.0588
$ws.Range("I89").Value = 3558.6667
$ws.Range("J89").Value = 2431.2
$ws.Range("K89").Value = 17793.3335
$ws.Range("L89").Value = 12156
$ws.Range("M89").Value = -12177.3335
$ws.Range("N89").Value = -23388
$ws.Range("H94").Value = 445.5
$ws.Range("I94").Value = 312.54544
$ws.Range("J94").Value = 933
$ws.Range("K94").Value = 312.54544
$ws.Range("L94").Value = 933
$ws.Range("M94").Value = 138.45456
$ws.Range("N94").Value = -1835
$ws.Range("H99").Value = 41127.08
$ws.Range("I99").Value = 85449.664
$ws.Range("J99").Value = 3136.2856
$ws.Range("K99").Value = 85449.664
$ws.Range("L99").Value = 3136.2856
$ws.Range("M99").Value = -83951.664
$ws.Range("N99").Value = -6132.2856
$ws.Range("H134").Value = 11653.863
$ws.Range("I134").Value = 11653.863
$ws.Range("K134").Value = 34961.589
$ws.Range("M134").Value = -32426.589
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 101.82143
$ws.Range("I7").Value = 73.181816
$ws.Range("K7").Value = 73.181816
$ws.Range("M7").Value = 39.818184
$ws.Range("H22").Value = 732.19354
$ws.Range("I22").Value = 736.3214
$ws.Range("J22").Value = 693.6667
$ws.Range("K22").Value = 736.3214
$ws.Range("L22").Value = 693.6667
$ws.Range("M22").Value = -386.3214
$ws.Range("N22").Value = -1393.6667
$ws.Range("H31").Value = 13046.315
$ws.Range("I31").Value = 1095.579
$ws.Range("J31").Value = 24997.053
$ws.Range("K31").Value = 1095.579
$ws.Range("L31").Value = 24997.053
$ws.Range("M31").Value = -800.579
$ws.Range("N31").Value = -25587.053
$ws.Range("H34").Value = 13046.315
$ws.Range("I34").Value = 1095.579
$ws.Range("J34").Value = 24997.053
$ws.Range("K34").Value = 1095.579
$ws.Range("L34").Value = 24997.053
$ws.Range("M34").Value = -893.579
$ws.Range("N34").Value = -25401.053
$ws.Range("H39").Value = 14823.333
$ws.Range("I39").Value = 14888
$ws.Range("K39").Value = 14888
$ws.Range("M39").Value = -14497
$ws.Range("H41").Value = 28958.309
$ws.Range("J41").Value = 35349.875
$ws.Range("L41").Value = 35349.875
$ws.Range("N41").Value = -36205.875
$ws.Range("H49").Value = 14823.333
$ws.Range("I49").Value = 14888
$ws.Range("K49").Value = 14888
$ws.Range("M49").Value = -14706
$ws.Range("H55").Value = 4891
$ws.Range("J55").Value = 4891
$ws.Range("L55").Value = 4891
$ws.Range("N55").Value = -5521
$ws.Range("H99").Value = 4158.6665
$ws.Range("I99").Value = 3970.2
$ws.Range("J99").Value = 4269.5293
$ws.Range("K99").Value = 3970.2
$ws.Range("L99").Value = 4269.5293
$ws.Range("M99").Value = -2472.2
$ws.Range("N99").Value = -7265.5293
$ws.Range("H126").Value = 4158.6665
$ws.Range("I126").Value = 3970.2
$ws.Range("J126").Value = 4269.5293
$ws.Range("K126").Value = 11910.6
$ws.Range("L126").Value = 12808.5879
$ws.Range("M126").Value = -9440.599999999999
$ws.Range("N126").Value = -17748.5879
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("H134").Value = 1320.7222
$ws.Range("I134").Value = 1320.7222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3962.1666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1427.1666
$ws.Range("N134").Value = $null
$ws.Range("H141").Value = 96938
$ws.Range("J141").Value = 96938
$ws.Range("L141").Value = 96938
$ws.Range("N141").Value = -107298

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 969.8333
$ws.Range("J23").Value = 994.2727
$ws.Range("L23").Value = 2982.8181
$ws.Range("N23").Value = -3452.8181
$ws.Range("H46").Value = 2354.1428
$ws.Range("J46").Value = 3533.3333
$ws.Range("L46").Value = 10599.9999
$ws.Range("N46").Value = -10781.9999
$ws.Range("H58").Value = 1200
$ws.Range("I58").Value = 1100
$ws.Range("J58").Value = 1600
$ws.Range("K58").Value = 3300
$ws.Range("L58").Value = 4800
$ws.Range("M58").Value = -3172
$ws.Range("N58").Value = -5056
$ws.Range("H98").Value = 1367.8889
$ws.Range("I98").Value = 498.66666
$ws.Range("J98").Value = 1802.5
$ws.Range("K98").Value = 1495.99998
$ws.Range("L98").Value = 5407.5
$ws.Range("M98").Value = 2.00001999999995
$ws.Range("N98").Value = -8403.5
$ws.Range("H107").Value = 684.8889
$ws.Range("J107").Value = 1111.75
$ws.Range("L107").Value = 3335.25
$ws.Range("N107").Value = -7175.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 28908.043
$ws.Range("J15").Value = 30325.934
$ws.Range("L15").Value = 30325.934
$ws.Range("N15").Value = -30901.934
$ws.Range("H81").Value = 28908.043
$ws.Range("J81").Value = 30325.934
$ws.Range("L81").Value = 30325.934
$ws.Range("N81").Value = -32321.934
$ws.Range("H84").Value = 28908.043
$ws.Range("J84").Value = 30325.934
$ws.Range("L84").Value = 90977.802
$ws.Range("N84").Value = -100961.802
$ws.Range("H97").Value = 1527.4667
$ws.Range("I97").Value = 400.3
$ws.Range("J97").Value = 3781.8
$ws.Range("K97").Value = 400.3
$ws.Range("L97").Value = 3781.8
$ws.Range("M97").Value = 95.69999999999999
$ws.Range("N97").Value = -4773.8
$ws.Range("H102").Value = 2418.6897
$ws.Range("I102").Value = 1374.6666
$ws.Range("K102").Value = 1374.6666
$ws.Range("M102").Value = 247.3334
$ws.Range("H113").Value = 3445.9375
$ws.Range("J113").Value = 4033.889
$ws.Range("L113").Value = 4033.889
$ws.Range("N113").Value = -8373.888999999999
$ws.Range("H122").Value = 4500.0625
$ws.Range("I122").Value = 7983.6665
$ws.Range("J122").Value = 2409.9
$ws.Range("K122").Value = 23950.9995
$ws.Range("L122").Value = 7229.700000000001
$ws.Range("M122").Value = -21500.9995
$ws.Range("N122").Value = -12129.7
$ws.Range("H126").Value = 3859.5833
$ws.Range("I126").Value = 3289.375
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9868.125
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -7398.125
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 6939.6313
$ws.Range("I132").Value = 7491
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 22473
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -19943
$ws.Range("N132").Value = -17057
$ws.Range("H136").Value = 15000
$ws.Range("J136").Value = 15000
$ws.Range("L136").Value = 45000
$ws.Range("N136").Value = -50100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = $null
$ws.Range("H22").Value = 18302.572
$ws.Range("I22").Value = 59619.5
$ws.Range("J22").Value = 1775.8
$ws.Range("K22").Value = 59619.5
$ws.Range("L22").Value = 1775.8
$ws.Range("M22").Value = -59324.5
$ws.Range("N22").Value = -2365.8
$ws.Range("H27").Value = 18302.572
$ws.Range("I27").Value = 59619.5
$ws.Range("J27").Value = 1775.8
$ws.Range("K27").Value = 59619.5
$ws.Range("L27").Value = 1775.8
$ws.Range("M27").Value = -59512.5
$ws.Range("N27").Value = -1989.8
$ws.Range("H31").Value = 2392.75
$ws.Range("I31").Value = 4372.5
$ws.Range("K31").Value = 4372.5
$ws.Range("M31").Value = -4124.5
$ws.Range("H32").Value = 6750.636
$ws.Range("J32").Value = 19999
$ws.Range("L32").Value = 19999
$ws.Range("N32").Value = -20633
$ws.Range("H35").Value = 2545.7827
$ws.Range("I35").Value = 1150.9
$ws.Range("K35").Value = 1150.9
$ws.Range("M35").Value = -814.9000000000001
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null
$ws.Range("H46").Value = 3986
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 4206.4287
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 4206.4287
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -4582.4287
$ws.Range("H51").Value = 41000
$ws.Range("J51").Value = 44000
$ws.Range("L51").Value = 44000
$ws.Range("N51").Value = -44956
$ws.Range("H63").Value = 70366.336
$ws.Range("I63").Value = 61038.5
$ws.Range("J63").Value = 73031.42999999999
$ws.Range("K63").Value = 61038.5
$ws.Range("L63").Value = 73031.42999999999
$ws.Range("M63").Value = -60289.5
$ws.Range("N63").Value = -74529.42999999999
$ws.Range("H66").Value = 70366.336
$ws.Range("I66").Value = 61038.5
$ws.Range("J66").Value = 73031.42999999999
$ws.Range("K66").Value = 183115.5
$ws.Range("L66").Value = 219094.29
$ws.Range("M66").Value = -179371.5
$ws.Range("N66").Value = -226582.29
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = $null
$ws.Range("H100").Value = 27750.75
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 52501.5
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 52501.5
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -53583.5
$ws.Range("H122").Value = 90333.664
$ws.Range("I122").Value = 6555.5557
$ws.Range("J122").Value = 341668
$ws.Range("K122").Value = 19666.6671
$ws.Range("L122").Value = 1025004
$ws.Range("M122").Value = -17216.6671
$ws.Range("N122").Value = -1029904
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = $null
$ws.Range("H136").Value = 18073.902
$ws.Range("I136").Value = 2852.3572
$ws.Range("J136").Value = 160141.67
$ws.Range("K136").Value = 8557.071599999999
$ws.Range("L136").Value = 480425.01
$ws.Range("M136").Value = -6007.071599999999
$ws.Range("N136").Value = -485525.01

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28129.133
$ws.Range("I62").Value = 4251
$ws.Range("J62").Value = 36812.09
$ws.Range("K62").Value = 4251
$ws.Range("L62").Value = 36812.09
$ws.Range("M62").Value = -3627
$ws.Range("N62").Value = -38060.09
$ws.Range("H65").Value = 28129.133
$ws.Range("I65").Value = 4251
$ws.Range("J65").Value = 36812.09
$ws.Range("K65").Value = 21255
$ws.Range("L65").Value = 184060.45
$ws.Range("M65").Value = -18135
$ws.Range("N65").Value = -190300.45
$ws.Range("H98").Value = 96648
$ws.Range("J98").Value = 96648
$ws.Range("L98").Value = 96648
$ws.Range("N98").Value = -102638
$ws.Range("H107").Value = 4221.5557
$ws.Range("I107").Value = 999.5
$ws.Range("J107").Value = 5142.143
$ws.Range("K107").Value = 2998.5
$ws.Range("L107").Value = 15426.429
$ws.Range("M107").Value = -1078.5
$ws.Range("N107").Value = -19266.429
$ws.Range("H122").Value = 2170.1428
$ws.Range("I122").Value = 2424
$ws.Range("J122").Value = 1831.6666
$ws.Range("K122").Value = 7272
$ws.Range("L122").Value = 5494.9998
$ws.Range("M122").Value = -4822
$ws.Range("N122").Value = -10394.9998
$ws.Range("H132").Value = 9035.3125
$ws.Range("I132").Value = 3004.889
$ws.Range("J132").Value = 16788.715
$ws.Range("K132").Value = 9014.667000000001
$ws.Range("L132").Value = 50366.145
$ws.Range("M132").Value = -6484.667000000001
$ws.Range("N132").Value = -55426.145

